$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The paragraph "Sum = sum+a[i]" is edited to become "sum = sum+mas[i]":
#   - "Sum" -> split into two runs "s" + "um" (same formatting)
#   - "a"   -> split into three runs "m" + "a" + "s" (same formatting)
#   - the "_GoBack" bookmark (normally tracking the last edit position)
#     is moved from its old location to right after the new "s" run.
# ------------------------------------------------------------------

# --- Step 1: "a" (in "+a[i]") becomes "m" + "a" + "s" ---
$text = $d.Content.Text
$idxA = $text.IndexOf("+a[i]") + 1   # position of the 'a' character

# Replace just the 'a' character with "mas"; because we replace the text of
# the existing 'a' run (rather than inserting next to it), the new text
# inherits that run's exact formatting (rFonts/color/sz/lang en-US).
$aRange = $d.Range($idxA, $idxA + 1)
$aRange.Text = "mas"

# Force Word to split "mas" into three separate runs (m | a | s) by briefly
# toggling a character formatting property on each sub-range; toggling it
# back to the original value keeps the visible formatting unchanged but
# causes the run to be split off from its neighbors.
$mRange = $d.Range($idxA, $idxA + 1)
$mRange.Bold = 1
$mRange.Bold = 0

$sRange = $d.Range($idxA + 2, $idxA + 3)
$sRange.Bold = 1
$sRange.Bold = 0

# --- Step 2: "Sum" becomes "s" + "um" ---
$text2 = $d.Content.Text
$idxSum = $text2.IndexOf("Sum")

$firstRange = $d.Range($idxSum, $idxSum + 1)
$firstRange.Text = "s"

$firstSplit = $d.Range($idxSum, $idxSum + 1)
$firstSplit.Bold = 1
$firstSplit.Bold = 0

# --- Step 3: move the "_GoBack" bookmark to right after the new "s" run ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newBookmarkPos = $idxA + 3   # right after "mas", before "["
$bookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
